$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.497.39"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.870.23"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  -1.30%  "
$ws.Range("D5").Value = "'315.15"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "'0.5070"
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("D8").Value = "'0.3902"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.08372"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'1.106"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "'41.77"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'6.222"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.875.91"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'20.46"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'7.282"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "'0.00001102"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'91.14"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'0.06735"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "'17.73"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'5.924"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "28.518.06"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "'11.09"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "'2.212"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "2.085.92"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'158.67"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").Value = "'20.61"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "'2.428"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "'127.12"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'0.1040"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "'1.042"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "'5.739"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").Value = "'3.624"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "'0.02456"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "'0.06580"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "'0.2166"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'8.900"
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("D39").Value = "'5.031"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "'1.180"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").Value = "'0.6374"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "'0.6012"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'13.10"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "'3.687"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'2.005"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "'1.215"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'122.49"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'0.06807"
$ws.Range("E51").Value = "  -0.79%  "
